$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "SUBSCRIPCIÓN" header in column G (7), same header row as the rest
$ws.Cells.Item(1, 7).Value = "SUBSCRIPCIÓN"

# Match the header style/fill of the existing header cells (copy F1's format onto G1)
$ws.Cells.Item(1, 6).Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)  # xlPasteFormats

# Column F (old last data column) now matches column E's width/outline grouping,
# and the new column G takes over the old "last column" width or 17
$ws.Columns.Item(6).ColumnWidth = 16.666666666666668   # -> stored width 17.5
$ws.Columns.Item(6).OutlineLevel = 1
$ws.Columns.Item(7).ColumnWidth = 16.166666666666668   # -> stored width 17

# Extend the autofilter to cover the new column
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:G1").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new autofilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$G`$1"
    }
}

# Reflect the new active cell/selection
$null = $ws.Range("G2").Select()
